$wb = $excel.ActiveWorkbook

# --- Add the new "Player Info" sheet, placed before "ODI Batting" ---
$battingSheet = $wb.Worksheets.Item("ODI Batting")
$playerInfo = $wb.Worksheets.Add($battingSheet)
$playerInfo.Name = "Player Info"

# Header row (bold, bordered, centered - matching the existing header style)
$playerInfo.Range("A1:D1").Font.Bold = $true
$playerInfo.Range("A1:D1").HorizontalAlignment = -4108
$playerInfo.Range("A1:D1").VerticalAlignment = -4160
$playerInfo.Range("A1:D1").Borders.LineStyle = 1

$playerInfo.Range("A1").Value = "ID"
$playerInfo.Range("B1").Value = "NAME"
$playerInfo.Range("C1").Value = "BATTING_HAND"
$playerInfo.Range("D1").Value = "BOWL_STYLE"

# Data row - keep the ID as text so it matches the original string typing
$playerInfo.Range("A2").NumberFormat = "@"
$playerInfo.Range("A2").Value = "3826"
$playerInfo.Range("B2").Value = "Cheteshwar Arvind Pujara"
$playerInfo.Range("C2").Value = "Right Handed"
$playerInfo.Range("D2").Value = "Right Arm Leg Break"

$playerInfo.Range("A1").Select()

# --- Update "ODI Batting" sheet: MATCH_CARD_LINK -> MATCH_CODE ---
# Re-fetch the worksheet by name since adding a new sheet shifted the
# earlier reference's focus.
$battingSheet = $wb.Worksheets.Item("ODI Batting")

$battingSheet.Range("D1").Value = "MATCH_CODE"

# Replace the full scorecard URLs with the bare match code, keeping the
# values stored as text (consistent with the rest of the sheet).
$battingSheet.Range("D2:D6").NumberFormat = "@"
$battingSheet.Range("D2").Value = "3544"
$battingSheet.Range("D3").Value = "3545"
$battingSheet.Range("D4").Value = "3643"
$battingSheet.Range("D5").Value = "3644"
$battingSheet.Range("D6").Value = "3645"
